$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 455, shifting existing rows 455:558 down to 456:559
$ws.Rows.Item(455).Insert()

# Populate the newly inserted row 455 with the new record's data.
# (Columns A, B, C, E, F, G, H, I, J, K, L, Q, R, T keep the same values
#  as the template row that used to occupy row 455.)
$ws.Range("A455").Value = 7
$ws.Range("B455").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C455").Value = "Ñuble"
$ws.Range("D455").Value = 44641
$ws.Range("E455").Value = 16
$ws.Range("F455").Value = "Fruta"
$ws.Range("G455").Value = 100108
$ws.Range("H455").Value = "Tropicales y subtropicales"
$ws.Range("I455").Value = 100108006
$ws.Range("J455").Value = "Plátano"
$ws.Range("K455").Value = "Sin especificar"
$ws.Range("L455").Value = "Pintón"
$ws.Range("M455").Value = 320
$ws.Range("N455").Value = 19000
$ws.Range("O455").Value = 21000
$ws.Range("P455").Value = 20125
$ws.Range("Q455").Value = "$/caja 20 kilos"
$ws.Range("R455").Value = "Ecuador"
$ws.Range("S455").Value = 1006
$ws.Range("T455").Value = 20
